$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: row number, optional A value, B value
$data = @(
  ,@(3, $null, 0.67298194597496086)
  ,@(4, $null, 0.67127144929829152)
  ,@(5, $null, 0.66953795732020382)
  ,@(6, $null, 0.66778166970734387)
  ,@(7, $null, 0.66600278840446681)
  ,@(8, $null, 0.66420151759574309)
  ,@(9, $null, 0.66237806366570795)
  ,@(10, $null, 0.66053263515986227)
  ,@(11, $null, 0.65866544274493821)
  ,@(12, 11, 0.6567766991688373)
  ,@(13, 12, 0.6548666192202538)
  ,@(14, 13, 0.65293541968799296)
  ,@(15, 14, 0.65130247056258039)
  ,@(16, 15, 0.64964707591482374)
  ,@(17, 16, 0.64796942673748459)
  ,@(18, 17, 0.64626971626513763)
  ,@(19, 18, 0.64454813993707893)
  ,@(20, 19, 0.64280489535987939)
  ,@(21, 20, 0.64104018226959214)
  ,@(22, 21, 0.63925420249362275)
  ,@(23, 22, 0.63744715991227396)
  ,@(24, 23, 0.6356192604199753)
  ,@(25, 24, 0.63377071188620704)
  ,@(26, 25, 0.63190172411612866)
  ,@(27, 26, 0.63032137889877915)
  ,@(28, 27, 0.62871931122030944)
  ,@(29, 28, 0.62709570592083108)
  ,@(30, 29, 0.62545075001004957)
  ,@(31, 30, 0.62378463263136963)
  ,@(32, 31, 0.62209754502565384)
  ,@(33, 32, 0.62038968049464771)
  ,@(34, 33, 0.61866123436407783)
  ,@(35, 34, 0.61691240394643698)
  ,@(36, 35, 0.61514338850346029)
  ,@(37, 36, 0.61335438920830998)
  ,@(38, 37, 0.61154560910747124)
  ,@(39, 38, 0.61001617321315404)
  ,@(40, 39, 0.60846571462620924)
  ,@(41, 40, 0.60689441223229912)
  ,@(42, 41, 0.6053024470167887)
  ,@(43, 42, 0.60369000203000767)
  ,@(44, 43, 0.60205726235217538)
  ,@(45, 44, 0.60040441505800446)
  ,@(46, 45, 0.59873164918098887)
  ,@(47, 46, 0.59703915567738608)
  ,@(48, 47, 0.59532712738990368)
  ,@(49, 48, 0.59359575901110107)
  ,@(50, 49, 0.59184524704651342)
  ,@(51, 50, 0.59036458848865836)
  ,@(52, 51, 0.58886309618547028)
  ,@(53, 52, 0.58734094707464068)
  ,@(54, 53, 0.585798320226953)
  ,@(55, 54, 0.58423539681117387)
  ,@(56, 55, 0.58265236005859711)
  ,@(57, 56, 0.58104939522724919)
  ,@(58, 57, 0.57942668956576471)
  ,@(59, 58, 0.57778443227694543)
  ,@(60, 59, 0.57612281448100766)
  ,@(61, 60, 0.57816016142613846)
)

foreach ($row in $data) {
    $r = $row[0]
    $a = $row[1]
    $b = $row[2]
    if ($a -ne $null) {
        $ws.Cells.Item($r, 1).Value = $a
    }
    $ws.Cells.Item($r, 2).Value = $b
}

$ws.Range("E10").Select()

